$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: new cell O5 = 1
$ws.Range("O5").Value = 1

# Row 7: G7 gets a value (keeps existing style s="11"); new cell N7 = 1
$ws.Range("G7").Value = 1
$ws.Range("N7").Value = 1

# Row 13: new cell O13 = 1
$ws.Range("O13").Value = 1

# Row 16: new cells N16 = 1, O16 = 1
$ws.Range("N16").Value = 1
$ws.Range("O16").Value = 1

# Row 17: new cell N17 = 1 (O17 already existed)
$ws.Range("N17").Value = 1

# Row 18: G18 gets a value (keeps existing style s="11"); new cell L18 = 1
$ws.Range("G18").Value = 1
$ws.Range("L18").Value = 1

# Row 21: G21 gets a value (keeps existing style s="11"); new cells L21 = 1, N21 = 1
$ws.Range("G21").Value = 1
$ws.Range("L21").Value = 1
$ws.Range("N21").Value = 1

# Row 25: new cell N25 = 1
$ws.Range("N25").Value = 1

# Row 27: new cell O27 = 1
$ws.Range("O27").Value = 1

# Update the active selection to L21 to match the saved view state
$ws.Range("L21").Select()
